$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.490.01"

$ws.Cells.Item(3, 4).Value = "3.057.67"
$ws.Cells.Item(3, 5).Value = "  -2.72%  "

$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.24%  "

$ws.Cells.Item(5, 4).Value = "590.05"
$ws.Cells.Item(5, 5).Value = "  -0.39%  "

$ws.Cells.Item(6, 4).Value = "154.79"
$ws.Cells.Item(6, 5).Value = "  +6.53%  "

$ws.Cells.Item(7, 5).Value = "  -0.18%  "

$ws.Cells.Item(8, 5).Value = "  +3.10%  "

$ws.Cells.Item(9, 4).Value = "3.074.48"
$ws.Cells.Item(9, 5).Value = "  -1.92%  "

$ws.Cells.Item(10, 5).Value = "  -1.99%  "

$ws.Cells.Item(11, 5).Value = "  -0.31%  "

$ws.Cells.Item(12, 5).Value = "  -0.08%  "

$ws.Cells.Item(13, 5).Value = "  +1.13%  "

$ws.Cells.Item(14, 5).Value = "  -2.18%  "

$ws.Cells.Item(15, 5).Value = "  -1.58%  "

$ws.Cells.Item(16, 4).Value = "3.564.43"
$ws.Cells.Item(16, 5).Value = "  -2.75%  "

$ws.Cells.Item(17, 4).Value = "7.22"
$ws.Cells.Item(17, 5).Value = "  -1.51%  "

$ws.Cells.Item(18, 4).Value = "63.426.63"

$ws.Cells.Item(19, 4).Value = "3.069.80"
$ws.Cells.Item(19, 5).Value = "  -2.27%  "

$ws.Cells.Item(20, 4).Value = "478.38"
$ws.Cells.Item(20, 5).Value = "  +2.12%  "

$ws.Cells.Item(21, 5).Value = "  +2.43%  "

$ws.Cells.Item(22, 5).Value = "  -2.07%  "

$ws.Cells.Item(23, 5).Value = "  +0.76%  "

$ws.Cells.Item(24, 5).Value = "  +3.26%  "

$ws.Cells.Item(25, 5).Value = "  -0.27%  "

$ws.Cells.Item(26, 4).Value = "81.25"
$ws.Cells.Item(26, 5).Value = "  -0.17%  "

$ws.Cells.Item(27, 4).Value = "10.03"
$ws.Cells.Item(27, 5).Value = "  +2.20%  "

$ws.Cells.Item(28, 5).Value = "  -0.23%  "

$ws.Cells.Item(29, 5).Value = "  -0.18%  "

$ws.Cells.Item(30, 5).Value = "  -0.60%  "

$ws.Cells.Item(31, 5).Value = "  -0.22%  "

$ws.Cells.Item(32, 5).Value = "  -1.89%  "

$ws.Cells.Item(34, 4).Value = "27.26"
$ws.Cells.Item(34, 5).Value = "  -1.91%  "

$ws.Cells.Item(35, 5).Value = "  +0.85%  "

$ws.Cells.Item(36, 5).Value = "  -1.49%  "

$ws.Cells.Item(37, 4).Value = "3.38"
$ws.Cells.Item(37, 5).Value = "  +4.81%  "

$ws.Cells.Item(38, 5).Value = "  -0.57%  "

$ws.Cells.Item(39, 4).Value = "2.23"
$ws.Cells.Item(39, 5).Value = "  -3.89%  "

$ws.Cells.Item(40, 4).Value = "9.39"
$ws.Cells.Item(40, 5).Value = "  +1.25%  "

$ws.Cells.Item(41, 4).Value = "50.62"
$ws.Cells.Item(41, 5).Value = "  -1.69%  "

$ws.Cells.Item(42, 4).Value = "445.85"
$ws.Cells.Item(42, 5).Value = "  -2.02%  "

$ws.Cells.Item(43, 5).Value = "  -2.69%  "

$ws.Cells.Item(44, 5).Value = "  -1.96%  "

$ws.Cells.Item(45, 2).Value = "Arweave"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(45, 4).Value = "40.18"
$ws.Cells.Item(45, 5).Value = "  +1.19%  "

$ws.Cells.Item(46, 2).Value = "Kaspa"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(46, 4).Value = "0.112"
$ws.Cells.Item(46, 5).Value = "  +3.57%  "

$ws.Cells.Item(47, 4).Value = "2.803.81"
$ws.Cells.Item(47, 5).Value = "  -3.68%  "

$ws.Cells.Item(48, 4).Value = "131.87"
$ws.Cells.Item(48, 5).Value = "  +1.24%  "

$ws.Cells.Item(49, 5).Value = "  +0.05%  "

$ws.Cells.Item(50, 4).Value = "25.26"
$ws.Cells.Item(50, 5).Value = "  +4.28%  "

$ws.Cells.Item(51, 5).Value = "  +0.95%  "
